$d = $word.ActiveDocument

# --- Fix typos / add missing accents in the last commentary paragraph ---
$d.Content.Find.Execute("progamas", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "programas", 2) | Out-Null

$d.Content.Find.Execute("desafio con", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "desafío con", 2) | Out-Null

$d.Content.Find.Execute("haciendonos", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "haciéndonos", 2) | Out-Null

$d.Content.Find.Execute("medio dificil", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "medio difícil", 2) | Out-Null

# --- Remove the trailing empty paragraph at the end of the document ---
# The very last paragraph mark of a Word document can't be deleted on its
# own (Word always keeps a final paragraph mark), so instead we delete the
# range spanning from just before the *previous* paragraph's mark through
# the end of the trailing empty paragraph - this merges the two paragraphs
# and effectively removes the empty one.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastParaText = $lastPara.Range.Text.Trim([char]13, [char]7)

if ($lastParaText -eq "") {
    $prevPara = $d.Paragraphs.Item($count - 1)
    $rng = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
    $rng.Delete()
}
